{"js": "// Fix the typo in the \"Practical Usage Condion\" heading so it reads\n// \"Practical Usage Condition\" (insert the missing \"ti\").\nconst body = context.document.body;\n\n// Search for the misspelled heading text. Using the distinctive misspelling\n// \"Condion\" keeps the search narrow and resilient to any surrounding\n// whitespace/run differences.\nconst results = body.search(\"Condion\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the misspelled word with the corrected spelling.\n  results.items[0].insertText(\"Condition\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Fix the typo in the \"Practical Usage Condion\" heading so it reads\n# \"Practical Usage Condition\" (insert the missing \"ti\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute(\"Condion\", $false, $false, $false, $false, $false, $true, 1, $false, \"Condition\", 2) | Out-Null\n"}
